$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift all Timestamp values (column A, rows 2-101) forward by 9 days.
# (Data originally covered 2025-11-22/23, corrected to 2025-12-01/02.)
for ($r = 2; $r -le 101; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value2 = $cell.Value2 + 9
}

# Fix the SGL production forecast values that were wrong for several
# intervals during the day (RO_RS deficit/excedent numbers).

# Rows 6-9: RO_RS deficit (E), Total Deficit (M) and System Direction (N)
# incorrectly carried a leftover 174.468 value -> should be 0.
foreach ($r in 6..9) {
    $ws.Cells.Item($r, 5).Value2 = 0    # E - Deficit_RO_RS
    $ws.Cells.Item($r, 13).Value2 = 0   # M - Total Deficit
    $ws.Cells.Item($r, 14).Value2 = 0   # N - System Direction
}

# Rows 18-21: same issue, leftover 181.327 value -> should be 0.
foreach ($r in 18..21) {
    $ws.Cells.Item($r, 5).Value2 = 0    # E - Deficit_RO_RS
    $ws.Cells.Item($r, 13).Value2 = 0   # M - Total Deficit
    $ws.Cells.Item($r, 14).Value2 = 0   # N - System Direction
}

# Rows 30-33: missing Excedent_RO_BG (B) / Deficit_RO_RS (E) values that
# should feed into Total Excedent (L) / Total Deficit (M) / System Direction (N).
foreach ($r in 30..33) {
    $ws.Cells.Item($r, 2).Value2 = 953.4999999999999   # B - Excedent_RO_BG
    $ws.Cells.Item($r, 5).Value2 = 243.9365             # E - Deficit_RO_RS
    $ws.Cells.Item($r, 12).Value2 = 953.4999999999999   # L - Total Excedent
    $ws.Cells.Item($r, 13).Value2 = 243.9365            # M - Total Deficit
    $ws.Cells.Item($r, 14).Value2 = 709.5634999999999   # N - System Direction
}

# Rows 34-37: missing Excedent_RO_BG (B) value feeding Total Excedent (L)
# and System Direction (N); Total Deficit (M) stays 0.
foreach ($r in 34..37) {
    $ws.Cells.Item($r, 2).Value2 = 1770.18   # B - Excedent_RO_BG
    $ws.Cells.Item($r, 12).Value2 = 1770.18  # L - Total Excedent
    $ws.Cells.Item($r, 14).Value2 = 1770.18  # N - System Direction
}

Write-Output "done"
